# Append two new data rows (39 and 40) to the bottom of Plan1's table,
# then update the sheet view's scroll position / selection to match the
# reviewer's final viewport.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$row39 = @(48, 16, 19, 14, 16, 8, 7, 47, 54, 32, 26, 71, 34, 22, 5)
$row40 = @(49, 15, 11, 13, 20, 11, 6, 47, 44, 14, 43, 76, 41, 21, 9)

for ($c = 0; $c -lt $row39.Length; $c++) {
    $ws.Cells.Item(39, $c + 1).Value = $row39[$c]
}

for ($c = 0; $c -lt $row40.Length; $c++) {
    $ws.Cells.Item(40, $c + 1).Value = $row40[$c]
}

# Move the viewport / selection to match the saved view state.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("J25").Select() | Out-Null
